$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace spaces with underscores in the category header names (row 1)
$ws.Range("B1").Value = "Smart_home_devices"
$ws.Range("C1").Value = "Smart_STBs"
$ws.Range("D1").Value = "Smart_TVs"
$ws.Range("E1").Value = "360_cameras"
$ws.Range("F1").Value = "VR_devices"
$ws.Range("H1").Value = "Other_portables"

# Move/update the active cell selection as recorded in the saved workbook
$ws.Range("J13").Select()
